$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 898.2857
$ws.Cells.Item(19, 9).Value = 950.25
$ws.Cells.Item(19, 10).Value = 829
$ws.Cells.Item(19, 11).Value = 950.25
$ws.Cells.Item(19, 12).Value = 829
$ws.Cells.Item(19, 13).Value = -775.25
$ws.Cells.Item(19, 14).Value = -1179
$ws.Cells.Item(118, 8).Value = 1594.5
$ws.Cells.Item(118, 9).Value = 792.6667
$ws.Cells.Item(118, 11).Value = 2378.0001
$ws.Cells.Item(118, 13).Value = -721.0001000000002
$ws.Cells.Item(132, 8).Value = 16081.786
$ws.Cells.Item(132, 9).Value = 16081.786
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 48245.358
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -45715.358
$ws.Cells.Item(132, 14).ClearContents()
$ws.Cells.Item(135, 8).Value = 1124.2858
$ws.Cells.Item(135, 9).Value = 1229.1666
$ws.Cells.Item(135, 10).Value = 495
$ws.Cells.Item(135, 11).Value = 11062.4994
$ws.Cells.Item(135, 12).Value = 4455
$ws.Cells.Item(135, 13).Value = -8527.499400000001
$ws.Cells.Item(135, 14).Value = -9525

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(26, 8).Value = 6833
$ws.Cells.Item(26, 9).Value = 6749.5
$ws.Cells.Item(26, 10).Value = 7000
$ws.Cells.Item(26, 11).Value = 6749.5
$ws.Cells.Item(26, 12).Value = 7000
$ws.Cells.Item(26, 13).Value = -6419.5
$ws.Cells.Item(26, 14).Value = -7660
$ws.Cells.Item(61, 8).Value = 4060.8572
$ws.Cells.Item(61, 10).Value = 6277.6665
$ws.Cells.Item(61, 12).Value = 6277.6665
$ws.Cells.Item(61, 14).Value = -6701.6665
$ws.Cells.Item(74, 8).Value = 3047
$ws.Cells.Item(74, 9).Value = 2686.6
$ws.Cells.Item(74, 11).Value = 2686.6
$ws.Cells.Item(74, 13).Value = -1812.6
$ws.Cells.Item(77, 8).Value = 3047
$ws.Cells.Item(77, 9).Value = 2686.6
$ws.Cells.Item(77, 11).Value = 13433
$ws.Cells.Item(77, 13).Value = -9065
$ws.Cells.Item(97, 8).Value = 734.5238000000001
$ws.Cells.Item(97, 9).Value = 719.2222
$ws.Cells.Item(97, 11).Value = 719.2222
$ws.Cells.Item(97, 13).Value = -223.2222
$ws.Cells.Item(122, 8).Value = 3040.6
$ws.Cells.Item(122, 9).Value = 2907.3333
$ws.Cells.Item(122, 11).Value = 8721.999899999999
$ws.Cells.Item(122, 13).Value = -6271.999899999999
$ws.Cells.Item(132, 8).Value = 3093.2307
$ws.Cells.Item(132, 9).Value = 3081.1
$ws.Cells.Item(132, 10).Value = 3133.6667
$ws.Cells.Item(132, 11).Value = 9243.299999999999
$ws.Cells.Item(132, 12).Value = 9401.000100000001
$ws.Cells.Item(132, 13).Value = -6713.299999999999
$ws.Cells.Item(132, 14).Value = -14461.0001
$ws.Cells.Item(136, 8).Value = 4060.8572
$ws.Cells.Item(136, 10).Value = 6277.6665
$ws.Cells.Item(136, 12).Value = 18832.9995
$ws.Cells.Item(136, 14).Value = -23932.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2049
$ws.Cells.Item(20, 9).Value = 1857
$ws.Cells.Item(20, 10).Value = 3009
$ws.Cells.Item(20, 11).Value = 1857
$ws.Cells.Item(20, 12).Value = 3009
$ws.Cells.Item(20, 13).Value = -1610
$ws.Cells.Item(20, 14).Value = -3503
$ws.Cells.Item(22, 8).Value = 1396.5714
$ws.Cells.Item(22, 9).Value = 1579.3334
$ws.Cells.Item(22, 10).Value = 300
$ws.Cells.Item(22, 11).Value = 1579.3334
$ws.Cells.Item(22, 12).Value = 300
$ws.Cells.Item(22, 13).Value = -1406.3334
$ws.Cells.Item(22, 14).Value = -646
$ws.Cells.Item(86, 8).Value = 6033.25
$ws.Cells.Item(86, 9).Value = 4561.1816
$ws.Cells.Item(86, 10).Value = 7832.4443
$ws.Cells.Item(86, 11).Value = 4561.1816
$ws.Cells.Item(86, 12).Value = 7832.4443
$ws.Cells.Item(86, 13).Value = -3438.1816
$ws.Cells.Item(86, 14).Value = -10078.4443
$ws.Cells.Item(89, 8).Value = 6033.25
$ws.Cells.Item(89, 9).Value = 4561.1816
$ws.Cells.Item(89, 10).Value = 7832.4443
$ws.Cells.Item(89, 11).Value = 22805.908
$ws.Cells.Item(89, 12).Value = 39162.2215
$ws.Cells.Item(89, 13).Value = -17189.908
$ws.Cells.Item(89, 14).Value = -50394.2215
$ws.Cells.Item(96, 8).Value = 39999.5
$ws.Cells.Item(96, 9).Value = 39999.5
$ws.Cells.Item(96, 11).Value = 39999.5
$ws.Cells.Item(96, 13).Value = -37253.5
$ws.Cells.Item(99, 8).Value = 2061
$ws.Cells.Item(99, 9).Value = 2252
$ws.Cells.Item(99, 11).Value = 2252
$ws.Cells.Item(99, 13).Value = -754
$ws.Cells.Item(134, 8).Value = 1557.2941
$ws.Cells.Item(134, 9).Value = 1498.3334
$ws.Cells.Item(134, 10).Value = 1999.5
$ws.Cells.Item(134, 11).Value = 4495.0002
$ws.Cells.Item(134, 12).Value = 5998.5
$ws.Cells.Item(134, 13).Value = -1960.0002
$ws.Cells.Item(134, 14).Value = -11068.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2605.4614
$ws.Cells.Item(58, 10).Value = 2808.4285
$ws.Cells.Item(58, 12).Value = 2808.4285
$ws.Cells.Item(58, 14).Value = -3214.4285
$ws.Cells.Item(105, 8).Value = 1986
$ws.Cells.Item(105, 9).Value = 1851.5
$ws.Cells.Item(105, 10).Value = 2299.8333
$ws.Cells.Item(105, 11).Value = 1851.5
$ws.Cells.Item(105, 12).Value = 2299.8333
$ws.Cells.Item(105, 13).Value = -104.5
$ws.Cells.Item(105, 14).Value = -5793.8333
$ws.Cells.Item(125, 8).Value = 4000
$ws.Cells.Item(125, 10).Value = 4000
$ws.Cells.Item(125, 12).Value = 4000
$ws.Cells.Item(125, 14).Value = -8920
$ws.Cells.Item(132, 8).Value = 4266.875
$ws.Cells.Item(132, 9).Value = 4077
$ws.Cells.Item(132, 10).Value = 4583.3335
$ws.Cells.Item(132, 11).Value = 12231
$ws.Cells.Item(132, 12).Value = 13750.0005
$ws.Cells.Item(132, 13).Value = -9701
$ws.Cells.Item(132, 14).Value = -18810.0005
$ws.Cells.Item(134, 8).Value = 1154.8572
$ws.Cells.Item(134, 9).Value = 1154.8572
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 11).Value = 3464.5716
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(134, 13).Value = -929.5715999999998
$ws.Cells.Item(134, 14).ClearContents()
$ws.Cells.Item(136, 8).Value = 2605.4614
$ws.Cells.Item(136, 10).Value = 2808.4285
$ws.Cells.Item(136, 12).Value = 8425.2855
$ws.Cells.Item(136, 14).Value = -13525.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 812.26666
$ws.Cells.Item(5, 9).Value = 788.2857
$ws.Cells.Item(5, 10).Value = 833.25
$ws.Cells.Item(5, 11).Value = 2364.8571
$ws.Cells.Item(5, 12).Value = 2499.75
$ws.Cells.Item(5, 13).Value = -2252.8571
$ws.Cells.Item(5, 14).Value = -2723.75
$ws.Cells.Item(14, 8).Value = 529.5454999999999
$ws.Cells.Item(14, 9).Value = 529.5454999999999
$ws.Cells.Item(14, 11).Value = 1588.6365
$ws.Cells.Item(14, 13).Value = -1415.6365
$ws.Cells.Item(15, 8).Value = 301
$ws.Cells.Item(15, 9).Value = 361.6
$ws.Cells.Item(15, 10).Value = 149.5
$ws.Cells.Item(15, 11).Value = 1084.8
$ws.Cells.Item(15, 12).Value = 448.5
$ws.Cells.Item(15, 13).Value = -944.8000000000002
$ws.Cells.Item(15, 14).Value = -728.5
$ws.Cells.Item(35, 8).Value = 800
$ws.Cells.Item(35, 9).Value = 800
$ws.Cells.Item(35, 11).Value = 2400
$ws.Cells.Item(35, 13).Value = -2112
$ws.Cells.Item(107, 8).Value = 330.66666
$ws.Cells.Item(107, 9).Value = 333.33334
$ws.Cells.Item(107, 10).Value = 325.33334
$ws.Cells.Item(107, 11).Value = 1000.00002
$ws.Cells.Item(107, 12).Value = 976.0000200000001
$ws.Cells.Item(107, 13).Value = 919.9999799999999
$ws.Cells.Item(107, 14).Value = -4816.00002
$ws.Cells.Item(108, 8).Value = 402.16666
$ws.Cells.Item(108, 9).Value = 402.16666
$ws.Cells.Item(108, 11).Value = 1206.49998
$ws.Cells.Item(108, 13).Value = 1673.50002
$ws.Cells.Item(135, 8).Value = 812.26666
$ws.Cells.Item(135, 9).Value = 788.2857
$ws.Cells.Item(135, 10).Value = 833.25
$ws.Cells.Item(135, 11).Value = 7094.571300000001
$ws.Cells.Item(135, 12).Value = 7499.25
$ws.Cells.Item(135, 13).Value = -4559.571300000001
$ws.Cells.Item(135, 14).Value = -12569.25
$ws.Cells.Item(140, 8).Value = 3774.9092
$ws.Cells.Item(140, 9).Value = 3190.5
$ws.Cells.Item(140, 11).Value = 9571.5
$ws.Cells.Item(140, 13).Value = -4391.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(99, 8).Value = 4500
$ws.Cells.Item(99, 9).Value = 4500
$ws.Cells.Item(99, 11).Value = 4500
$ws.Cells.Item(99, 13).Value = -2254
$ws.Cells.Item(107, 8).Value = 1369.4286
$ws.Cells.Item(107, 9).Value = 1597.4286
$ws.Cells.Item(107, 10).Value = 1141.4286
$ws.Cells.Item(107, 11).Value = 1597.4286
$ws.Cells.Item(107, 12).Value = 1141.4286
$ws.Cells.Item(107, 13).Value = 322.5714
$ws.Cells.Item(107, 14).Value = -4981.4286
$ws.Cells.Item(113, 8).Value = 1000
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 13).ClearContents()
$ws.Cells.Item(124, 8).Value = 0
$ws.Cells.Item(124, 10).Value = 0
$ws.Cells.Item(124, 12).Value = 0
$ws.Cells.Item(124, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 4057.6206
$ws.Cells.Item(132, 9).Value = 2946.2
$ws.Cells.Item(132, 10).Value = 6527.4443
$ws.Cells.Item(132, 11).Value = 8838.599999999999
$ws.Cells.Item(132, 12).Value = 19582.3329
$ws.Cells.Item(132, 13).Value = -6308.599999999999
$ws.Cells.Item(132, 14).Value = -24642.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 2597.35
$ws.Cells.Item(61, 9).Value = 1246.75
$ws.Cells.Item(61, 11).Value = 1246.75
$ws.Cells.Item(61, 13).Value = -1044.75
$ws.Cells.Item(113, 8).Value = 2597.35
$ws.Cells.Item(113, 9).Value = 1246.75
$ws.Cells.Item(113, 11).Value = 1246.75
$ws.Cells.Item(113, 13).Value = 923.25
$ws.Cells.Item(125, 8).Value = 228135
$ws.Cells.Item(125, 10).Value = 228135
$ws.Cells.Item(125, 12).Value = 228135
$ws.Cells.Item(125, 14).Value = -237975

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(103, 8).Value = 35936.375
$ws.Cells.Item(103, 10).Value = 35936.375
$ws.Cells.Item(103, 12).Value = 35936.375
$ws.Cells.Item(103, 14).Value = -38280.375
$ws.Cells.Item(136, 8).Value = 2911.1365
$ws.Cells.Item(136, 9).Value = 2190.3125
$ws.Cells.Item(136, 11).Value = 6570.9375
$ws.Cells.Item(136, 13).Value = -4020.9375

